$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top of the data (row 2), shifting existing rows down
$ws.Rows("2:2").Insert()

# Fill in the new "2023" entry
$ws.Range("A2").Value = "tendencias del campo de la comunicación estratégica"
$ws.Range("B2").Value = 2023
$ws.Range("C2").Value = "Universidad Central"
$ws.Range("D2").Value = "Bogotá, Colombia"

# Re-fit column widths to the new content, as Excel does automatically
# (equivalent to Excel's own AutoFit recalculation for the new data)
$ws.Columns("A").ColumnWidth = 76.66666666666667
$ws.Columns("B").ColumnWidth = 7.5
$ws.Columns("C").ColumnWidth = 60.0
$ws.Columns("D").ColumnWidth = 23.666666666666668
$ws.Columns("E").ColumnWidth = 6.166666666666667

# Update the active selection to match the saved state
$ws.Range("A11").Select() | Out-Null
